$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.751.65"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.890.77"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "246.73"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2923"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "0.06521"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "22.47"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "0.7411"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.889.03"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "96.77"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "5.244"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "284.98"
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").Value = "30.744.17"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "13.27"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007522"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.138.46"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "5.310"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "0.9995"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "6.264"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "9.233"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "164.41"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").Value = "18.96"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "1.919"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "0.09774"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").Value = "1.488"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "4.301"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "4.195"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").Value = "0.6992"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "0.01900"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "2.839"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").Value = "76.02"
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("D41").Value = "6.285"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "2.011"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "0.4282"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "0.8336"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "101.80"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "9.561"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").Value = "7.024"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").Value = "35.49"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "909.76"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "0.05781"
$ws.Range("E51").Value = "  +2.59%  "
